$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.120.49"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.654.49"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.63"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2605"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06344"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07798"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.679.58"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.498"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5467"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "0.0₅8193"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.30"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "26.120.34"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.575"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.30"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.05"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.028"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.09"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1240"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.235"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05899"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.522"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.244"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.588"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9524"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.786"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5680"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01616"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.828"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8500"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "1.031.79"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.01"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").Value = "1.800.51"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.14"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4301"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.885"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05166"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.472"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09696"
$ws.Range("E51").Value = "  -0.07%  "
